$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.995.51'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.339.06'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.45'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.31'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  -4.28%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.84'
$ws.Range('E10').Value = '  -4.24%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '52.42'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.90'
$ws.Range('E15').Value = '  +5.32%  '
$ws.Range('D16').Value = '2.342.34'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.831'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '42.921.45'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0910'
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.72'
$ws.Range('E20').Value = '  -5.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '69.19'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.80'
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.01'
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.65'
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('E29').Value = '  +1.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.24'
$ws.Range('E30').Value = '  -4.33%  '
$ws.Range('E31').Value = '  -4.08%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '161.89'
$ws.Range('E32').Value = '  -3.96%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  -3.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.66'
$ws.Range('E35').Value = '  +4.45%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.38'
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('E40').Value = '  -5.17%  '
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('E42').Value = '  -2.63%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.57'
$ws.Range('E43').Value = '  +3.07%  '
$ws.Range('D44').Value = '2.020.62'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0284'
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '18.89'
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.25'
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.93'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '55.78'
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.90'
$ws.Range('E50').Value = '  -2.14%  '
$ws.Range('D51').Value = '2.564.43'
$ws.Range('E51').Value = '  +1.08%  '
